# Novos produtos no index
# Update the "TIPO" (column D) tags for a set of rows so that certain
# products are flagged as "Index" or "Mais Vendidos" entries, and adjust
# the sheet view / column D width to match the author's final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Column D (TIPO) content updates -------------------------------------
$ws.Range("D2").Value  = "'Pulseira, Mais Vendidos',"
$ws.Range("D3").Value  = "'Pulseira, Index',"
$ws.Range("D4").Value  = "'Pulseira, Index',"
$ws.Range("D6").Value  = "'Pulseira, Index',"
$ws.Range("D7").Value  = "'Colar, Index',"
$ws.Range("D8").Value  = "'Conjunto, Mais Vendidos',"
$ws.Range("D11").Value = "'Colar, Index',"
$ws.Range("D12").Value = "'Pulseira, Index',"
$ws.Range("D13").Value = "'Colar, Index',"
$ws.Range("D17").Value = "'Colar, Index',"
$ws.Range("D21").Value = "'Brinco, Mais Vendidos',"
$ws.Range("D25").Value = "'Brinco, Index',"
$ws.Range("D28").Value = "'Brinco, Index',"
$ws.Range("D31").Value = "'Brinco, Index',"
$ws.Range("D35").Value = "'Brinco, Index',"
$ws.Range("D36").Value = "'Brinco, Mais Vendidos',"

# --- Column D width (now needs to fit the longer tag text) ----------------
$ws.Columns.Item(4).ColumnWidth = 22.6640625

# --- Sheet view / selection changes ---------------------------------------
$ws.Range("D5").Select()
$win = $excel.ActiveWindow
$win.Zoom = 120
$win.ScrollColumn = 2
